# increase maximum role columns in uploader to 12
#
# Adds 6 more "Role" blocks of 5 columns each (Role8 duplicated, then
# Role9..Role12) to the header row of the "Specialized Individual Template"
# and "Base Template" worksheets, mirroring the existing Role1..Role7
# pattern (columns L..AT). Also updates the active-sheet/selection state
# to match the post-edit workbook.

$wb = $excel.ActiveWorkbook

# New columns AU..BX (30 columns): Role8 (x2, duplicated in source), Role9,
# Role10, Role11, Role12 - each a 5-column block of ID/Name/Description/Type/Value.
$newCols = @("AU", "AV", "AW", "AX", "AY", "AZ", "BA", "BB", "BC", "BD", "BE", "BF", "BG", "BH", "BI", "BJ", "BK", "BL", "BM", "BN", "BO", "BP", "BQ", "BR", "BS", "BT", "BU", "BV", "BW", "BX")
$newLabels = @("Role8 ID", "Role8 Name", "Role8 Description", "Role8 Type", "Role8 Value", "Role8 ID", "Role8 Name", "Role8 Description", "Role8 Type", "Role8 Value", "Role9 ID", "Role9 Name", "Role9 Description", "Role9 Type", "Role9 Value", "Role10 ID", "Role10 Name", "Role10 Description", "Role10 Type", "Role10 Value", "Role11 ID", "Role11 Name", "Role11 Description", "Role11 Type", "Role11 Value", "Role12 ID", "Role12 Name", "Role12 Description", "Role12 Type", "Role12 Value")
# Source columns (existing Role7 block: L..AT repeats every 5 cols) whose
# formatting (fill/font/border style) each new column should copy - this
# keeps the banded "ID/Name/Description/Type/Value" style (s="1" four times
# then s="2") without minting new style entries.
$srcCols = @("AP", "AQ", "AR", "AS", "AT", "AP", "AQ", "AR", "AS", "AT", "AP", "AQ", "AR", "AS", "AT", "AP", "AQ", "AR", "AS", "AT", "AP", "AQ", "AR", "AS", "AT", "AP", "AQ", "AR", "AS", "AT")

function Add-RoleColumns($ws) {
    for ($i = 0; $i -lt $newCols.Count; $i++) {
        $destCell = $ws.Range($newCols[$i] + "1")
        $destCell.Value = $newLabels[$i]
        $ws.Range($srcCols[$i] + "1").Copy()
        $destCell.PasteSpecial(-4122)  # xlPasteFormats
    }
}

# --- Sheet 4: "Specialized Individual Template" ---
$ws4 = $wb.Worksheets.Item(4)
Add-RoleColumns $ws4

# --- Sheet 6: "Base Template" ---
$ws6 = $wb.Worksheets.Item(6)
Add-RoleColumns $ws6

# Selection / active-sheet bookkeeping to mirror the recorded edit: the user
# ended up on sheet 4 with L1:BY1 selected, after having last clicked
# through sheet 6 (topLeftCell BL1, selection L1:BX1 anchored at BX1).
$ws6.Activate()
$ws6.Range("BX1").Activate()
$ws6.Range("L1:BX1").Select()

$ws4.Activate()
$ws4.Range("L1:BY1").Select()
